$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Value = 1.44
$ws.Range("Q2").Value = 2.63
$ws.Range("P3").Value = 1.33
$ws.Range("AA4").Value = 7.5
$ws.Range("AB4").Value = 13
$ws.Range("AC4").Value = 41
$ws.Range("AJ4").Value = 21
$ws.Range("J4").Value = 1.03
$ws.Range("L4").Value = 1.18
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 2.25
$ws.Range("J5").Value = 1.05
$ws.Range("L5").Value = 1.3
$ws.Range("N5").Value = 1.98
$ws.Range("O5").Value = 1.83
$ws.Range("J6").Value = 1.13
$ws.Range("L6").Value = 1.62
$ws.Range("J25").Value = 1.02
$ws.Range("L25").Value = 1.13
$ws.Range("O25").Value = 2.35
$ws.Range("AA27").Value = 5.8
$ws.Range("AB27").Value = 15.5
$ws.Range("AC27").Value = 80
$ws.Range("AE27").Value = 4.65
$ws.Range("AF27").Value = 5.8
$ws.Range("AH27").Value = 10
$ws.Range("G27").Value = 4.9
$ws.Range("H27").Value = 3.3
$ws.Range("P27").Value = 1.47
$ws.Range("Q27").Value = 2.22
$ws.Range("T27").Value = 9.25
$ws.Range("U27").Value = 21
$ws.Range("V27").Value = 13.5
$ws.Range("W27").Value = 70
$ws.Range("Y27").Value = 50
$ws.Range("Z27").Value = 7.5
$ws.Range("AB28").Value = 11.75
$ws.Range("AE28").Value = 6
$ws.Range("AF28").Value = 7.9
$ws.Range("AG28").Value = 7.2
$ws.Range("AH28").Value = 14.5
$ws.Range("AI28").Value = 13.5
$ws.Range("G28").Value = 3.45
$ws.Range("H28").Value = 3.15
$ws.Range("I28").Value = 1.98
$ws.Range("T28").Value = 8.5
$ws.Range("U28").Value = 15
$ws.Range("V28").Value = 10
$ws.Range("W28").Value = 37
$ws.Range("X28").Value = 25
$ws.Range("AA32").Value = 11.25
$ws.Range("AB32").Value = 32
$ws.Range("AC32").Value = 175
$ws.Range("AE32").Value = 6.6
$ws.Range("AF32").Value = 5.5
$ws.Range("AG32").Value = 9.5
$ws.Range("AH32").Value = 6.7
$ws.Range("AI32").Value = 11.5
$ws.Range("AJ32").Value = 37
$ws.Range("G32").Value = 11.75
$ws.Range("H32").Value = 5.4
$ws.Range("I32").Value = 1.22
$ws.Range("L32").Value = 1.2
$ws.Range("M32").Value = 3.6
$ws.Range("N32").Value = 1.62
$ws.Range("O32").Value = 2.05
$ws.Range("R32").Value = 2.25
$ws.Range("S32").Value = 1.5
$ws.Range("T32").Value = 27
$ws.Range("U32").Value = 90
$ws.Range("V32").Value = 40
$ws.Range("W32").Value = 450
$ws.Range("X32").Value = 200
$ws.Range("Y32").Value = 175
$ws.Range("Z32").Value = 12
$ws.Range("AB34").Value = 13
$ws.Range("AD34").Value = 126
$ws.Range("G34").Value = 1.73
$ws.Range("H34").Value = 4
$ws.Range("I34").Value = 4.2
$ws.Range("R34").Value = 1.57
$ws.Range("S34").Value = 2.25
$ws.Range("T34").Value = 9.5
$ws.Range("U34").Value = 10
$ws.Range("Z34").Value = 17
$ws.Range("N35").Value = 1.6
$ws.Range("N36").Value = 1.95
$ws.Range("O36").Value = 1.85
$ws.Range("N39").Value = 2.88
$ws.Range("O39").Value = 1.4
$ws.Range("AC42").Value = 67
$ws.Range("AD42").Value = 450
$ws.Range("AG42").Value = 34
$ws.Range("G42").Value = 1.13
$ws.Range("N42").Value = 1.22
$ws.Range("O42").Value = 4
$ws.Range("J43").Value = 1.03
$ws.Range("K43").Value = 10.5
$ws.Range("AA44").Value = 6.3
$ws.Range("AE44").Value = 8.5
$ws.Range("AF44").Value = 12
$ws.Range("AG44").Value = 9
$ws.Range("AH44").Value = 24
$ws.Range("AI44").Value = 18.5
$ws.Range("AJ44").Value = 26
$ws.Range("G44").Value = 2.87
$ws.Range("H44").Value = 3.2
$ws.Range("I44").Value = 2.32
$ws.Range("L44").Value = 1.26
$ws.Range("M44").Value = 3.15
$ws.Range("N44").Value = 1.78
$ws.Range("O44").Value = 1.82
$ws.Range("S44").Value = 2.05
$ws.Range("T44").Value = 9.75
$ws.Range("U44").Value = 15.5
$ws.Range("V44").Value = 10.25
$ws.Range("W44").Value = 35
$ws.Range("X44").Value = 23
$ws.Range("Y44").Value = 29
$ws.Range("Z44").Value = 10.25
$ws.Range("AE45").Value = 10.25
$ws.Range("AF45").Value = 14.5
$ws.Range("AG45").Value = 10
$ws.Range("AH45").Value = 30
$ws.Range("AI45").Value = 20
$ws.Range("AJ45").Value = 27
$ws.Range("G45").Value = 2.3
$ws.Range("I45").Value = 2.65
$ws.Range("R45").Value = 1.6
$ws.Range("S45").Value = 2.07
$ws.Range("T45").Value = 9.5
$ws.Range("U45").Value = 12.5
$ws.Range("V45").Value = 9.25
$ws.Range("W45").Value = 23
$ws.Range("X45").Value = 17.5
$ws.Range("Y45").Value = 25
$ws.Range("Z45").Value = 12.5
$ws.Range("AA48").Value = 7.2
$ws.Range("AB48").Value = 15.5
$ws.Range("AD48").Value = 600
$ws.Range("AE48").Value = 12
$ws.Range("AF48").Value = 24
$ws.Range("AG48").Value = 14
$ws.Range("AI48").Value = 40
$ws.Range("G48").Value = 1.72
$ws.Range("H48").Value = 3.7
$ws.Range("I48").Value = 4.15
$ws.Range("L48").Value = 1.26
$ws.Range("M48").Value = 3.15
$ws.Range("O48").Value = 1.83
$ws.Range("R48").Value = 1.75
$ws.Range("S48").Value = 1.85
$ws.Range("T48").Value = 7.1
$ws.Range("U48").Value = 8.25
$ws.Range("W48").Value = 13.5
$ws.Range("X48").Value = 13.5
$ws.Range("Y48").Value = 26
$ws.Range("Z48").Value = 10.75
$ws.Range("AA51").Value = 6.5
$ws.Range("AB51").Value = 13.5
$ws.Range("AC51").Value = 60
$ws.Range("AD51").Value = 450
$ws.Range("AF51").Value = 21
$ws.Range("AG51").Value = 12.5
$ws.Range("AH51").Value = 60
$ws.Range("AI51").Value = 35
$ws.Range("AJ51").Value = 40
$ws.Range("G51").Value = 1.91
$ws.Range("H51").Value = 3.3
$ws.Range("I51").Value = 3.75
$ws.Range("L51").Value = 1.27
$ws.Range("M51").Value = 3.1
$ws.Range("N51").Value = 1.8
$ws.Range("O51").Value = 1.8
$ws.Range("R51").Value = 1.65
$ws.Range("S51").Value = 1.98
$ws.Range("T51").Value = 7.7
$ws.Range("U51").Value = 9.75
$ws.Range("W51").Value = 17.5
$ws.Range("X51").Value = 14.5
$ws.Range("Y51").Value = 24
$ws.Range("Z51").Value = 10
$ws.Range("AB53").Value = 13
$ws.Range("AC53").Value = 55
$ws.Range("AD53").Value = 400
$ws.Range("AE53").Value = 10.25
$ws.Range("AF53").Value = 17.5
$ws.Range("AG53").Value = 11.25
$ws.Range("AI53").Value = 28
$ws.Range("AJ53").Value = 32
$ws.Range("M53").Value = 3.35
$ws.Range("N53").Value = 1.83
$ws.Range("O53").Value = 1.87
$ws.Range("Q53").Value = 2.7
$ws.Range("T53").Value = 8.25
$ws.Range("U53").Value = 11
$ws.Range("X53").Value = 16
$ws.Range("Y53").Value = 24
$ws.Range("AA62").Value = 7.2
$ws.Range("AB62").Value = 12
$ws.Range("AE62").Value = 10.5
$ws.Range("AF62").Value = 13.5
$ws.Range("AG62").Value = 9.5
$ws.Range("AH62").Value = 26
$ws.Range("AI62").Value = 18
$ws.Range("AJ62").Value = 23
$ws.Range("G62").Value = 2.57
$ws.Range("H62").Value = 3.65
$ws.Range("I62").Value = 2.42
$ws.Range("K62").Value = 8.5
$ws.Range("L62").Value = 1.21
$ws.Range("M62").Value = 3.95
$ws.Range("R62").Value = 1.55
$ws.Range("S62").Value = 2.3
$ws.Range("T62").Value = 11
$ws.Range("U62").Value = 15
$ws.Range("V62").Value = 9.75
$ws.Range("W62").Value = 29
$ws.Range("X62").Value = 19
$ws.Range("Z62").Value = 8.5
